$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2248520710059172
$ws.Range("C2").Value = 0.5118343195266272
$ws.Range("J2").Value = 0.02071005917159763
$ws.Range("P2").Value = 0.1597633136094675
$ws.Range("S2").Value = 0.08284023668639054
$ws.Range("B3").Value = 0.005617977528089887
$ws.Range("C3").Value = 0.01685393258426966
$ws.Range("J3").Value = 0.0449438202247191
$ws.Range("P3").Value = 0.6910112359550562
$ws.Range("S3").Value = 0.2415730337078652
$ws.Range("J4").Value = 0.1
$ws.Range("P4").Value = 0.575
$ws.Range("S4").Value = 0.325
$ws.Range("B6").Value = 0.05128205128205128
$ws.Range("D6").Value = 0.008547008547008548
$ws.Range("F6").Value = 0.07692307692307693
$ws.Range("J6").Value = 0.2991452991452991
$ws.Range("O6").Value = 0.03846153846153846
$ws.Range("Q6").Value = 0.1452991452991453
$ws.Range("R6").Value = 0.05555555555555555
$ws.Range("S6").Value = 0.3247863247863248
$ws.Range("B7").Value = 0.0797872340425532
$ws.Range("D7").Value = 0.01595744680851064
$ws.Range("E7").Value = 0.01063829787234043
$ws.Range("F7").Value = 0.04787234042553191
$ws.Range("J7").Value = 0.1542553191489362
$ws.Range("O7").Value = 0.02659574468085106
$ws.Range("Q7").Value = 0.148936170212766
$ws.Range("R7").Value = 0.0851063829787234
$ws.Range("S7").Value = 0.4308510638297872
$ws.Range("B8").Value = 0.1086519114688129
$ws.Range("D8").Value = 0.01207243460764587
$ws.Range("F8").Value = 0.05835010060362173
$ws.Range("J8").Value = 0.1106639839034205
$ws.Range("O8").Value = 0.03018108651911469
$ws.Range("Q8").Value = 0.1629778672032193
$ws.Range("R8").Value = 0.07243460764587525
$ws.Range("S8").Value = 0.4446680080482898
$ws.Range("B9").Value = 0.06132075471698113
$ws.Range("D9").Value = 0.01415094339622642
$ws.Range("E9").Value = 0.004716981132075472
$ws.Range("F9").Value = 0.08962264150943396
$ws.Range("J9").Value = 0.1462264150943396
$ws.Range("O9").Value = 0.009433962264150943
$ws.Range("Q9").Value = 0.2169811320754717
$ws.Range("R9").Value = 0.08490566037735849
$ws.Range("S9").Value = 0.3726415094339622
$ws.Range("B10").Value = 0.1267605633802817
$ws.Range("D10").Value = 0.02190923317683881
$ws.Range("F10").Value = 0.05320813771517997
$ws.Range("J10").Value = 0.1267605633802817
$ws.Range("O10").Value = 0.03051643192488263
$ws.Range("Q10").Value = 0.1846635367762128
$ws.Range("R10").Value = 0.06729264475743349
$ws.Range("S10").Value = 0.3888888888888889
$ws.Range("G11").Value = 0.1807580174927114
$ws.Range("J11").Value = 0.1020408163265306
$ws.Range("K11").Value = 0.2361516034985423
$ws.Range("L11").Value = 0.4402332361516035
$ws.Range("S11").Value = 0.04081632653061224
$ws.Range("G12").Value = 0.6883116883116883
$ws.Range("J12").Value = 0.2207792207792208
$ws.Range("K12").Value = 0.006493506493506494
$ws.Range("L12").Value = 0.01948051948051948
$ws.Range("S12").Value = 0.06493506493506493
$ws.Range("F13").Value = 0.02380952380952381
$ws.Range("G13").Value = 0.5238095238095238
$ws.Range("J13").Value = 0.3571428571428572
$ws.Range("S13").Value = 0.09523809523809523
$ws.Range("F15").Value = 0.02880658436213992
$ws.Range("H15").Value = 0.1604938271604938
$ws.Range("I15").Value = 0.09053497942386832
$ws.Range("J15").Value = 0.3127572016460906
$ws.Range("K15").Value = 0.06584362139917696
$ws.Range("M15").Value = 0.01234567901234568
$ws.Range("O15").Value = 0.03703703703703703
$ws.Range("S15").Value = 0.2921810699588477
$ws.Range("F16").Value = 0.03092783505154639
$ws.Range("H16").Value = 0.1752577319587629
$ws.Range("I16").Value = 0.09278350515463918
$ws.Range("J16").Value = 0.3556701030927835
$ws.Range("K16").Value = 0.1030927835051546
$ws.Range("M16").Value = 0.02577319587628866
$ws.Range("O16").Value = 0.06701030927835051
$ws.Range("S16").Value = 0.1494845360824742
$ws.Range("F17").Value = 0.03087885985748218
$ws.Range("H17").Value = 0.1876484560570071
$ws.Range("I17").Value = 0.1163895486935867
$ws.Range("J17").Value = 0.3752969121140142
$ws.Range("K17").Value = 0.07600950118764846
$ws.Range("M17").Value = 0.007125890736342043
$ws.Range("O17").Value = 0.05938242280285035
$ws.Range("S17").Value = 0.1472684085510689
$ws.Range("F18").Value = 0.02312138728323699
$ws.Range("H18").Value = 0.1098265895953757
$ws.Range("I18").Value = 0.1040462427745665
$ws.Range("J18").Value = 0.4450867052023121
$ws.Range("K18").Value = 0.09248554913294797
$ws.Range("M18").Value = 0.02890173410404624
$ws.Range("O18").Value = 0.05202312138728324
$ws.Range("S18").Value = 0.1445086705202312
$ws.Range("F19").Value = 0.02162941600576784
$ws.Range("H19").Value = 0.2386445565969719
$ws.Range("I19").Value = 0.07570295602018745
$ws.Range("J19").Value = 0.3446286950252343
$ws.Range("K19").Value = 0.1196827685652487
$ws.Range("M19").Value = 0.01946647440519106
$ws.Range("N19").Value = 0.00144196106705119
$ws.Range("O19").Value = 0.06488824801730353
$ws.Range("S19").Value = 0.113914924297044
